$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.992433488368988
$ws.Range("B1").Value = 1.915817856788635
$ws.Range("C1").Value = 3.435439348220825
$ws.Range("D1").Value = 3.155739068984985
$ws.Range("E1").Value = 0.3562227487564087
